$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NEW")

$row = 73

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell $ws.Cells.Item($row, 1) "7229"
Set-TextCell $ws.Cells.Item($row, 2) "9/16/2025"
$ws.Cells.Item($row, 3).Value = "AZURDUY JUANA 2627"
Set-TextCell $ws.Cells.Item($row, 4) "13"
$ws.Cells.Item($row, 5).Value = "Pendiente ADM"
$ws.Cells.Item($row, 6).Value = "NEW"
$ws.Cells.Item($row, 7).Value = "Pendiente"
$ws.Cells.Item($row, 8).Value = "Colocar columna para pedir traspaso de nodo propio"
$ws.Cells.Item($row, 9).Value = 1
$ws.Cells.Item($row, 10).Value = "Cambio"
$ws.Cells.Item($row, 11).Value = "Nodo Teco"
$ws.Cells.Item($row, 12).Value = "Pasante"
$ws.Cells.Item($row, 13).Value = -58.469008
$ws.Cells.Item($row, 14).Value = -34.552083
$ws.Cells.Item($row, 15).Value = "Saavedra"
$ws.Cells.Item($row, 16).Value = "Capital Norte"
